$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 6: reuse the existing (previously-unused) "SingleUseId5" string for the TEXT ID column.
$ws.Range("B6").Value = "SingleUseId5"

# Row 7: introduce a brand-new unique text id "SingleUseId6".
$ws.Range("B7").Value = "SingleUseId6"

# New row 8: "Device settings:" entry, using the same Typography/Alignment/Direction as other
# left-aligned rows (row 4/5 style: Typography_00, Left, LTR).
$ws.Range("B8").Value = "SingleUseId7"
$ws.Range("C8").Value = "Typography_00"
$ws.Range("D8").Value = "Left"
$ws.Range("E8").Value = "LTR"
$ws.Range("F8").Value = "Device settings:"

# New row 9: "Work mode:" entry, new screen text.
$ws.Range("B9").Value = "SingleUseId8"
$ws.Range("C9").Value = "Typography_00"
$ws.Range("D9").Value = "Left"
$ws.Range("E9").Value = "LTR"
$ws.Range("F9").Value = "Work mode:"
